$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 193, shifting the existing rows 193-239 down to 196-242.
$ws.Rows("193:195").Insert()

# --- Common columns shared by all three new rows ---
$ws.Range("A193:A195").Value = 9
$ws.Range("B193:B195").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C193:C195").Value = "Metropolitana"
$ws.Range("D193:D195").Value = 44551
$ws.Range("E193:E195").Value = 13
$ws.Range("F193:F195").Value = "Fruta"
$ws.Range("G193:G195").Value = 100103
$ws.Range("H193:H195").Value = "Frutos de hueso (carozo)"
$ws.Range("I193:I195").Value = 100103001
$ws.Range("J193:J195").Value = "Cereza"

# --- Row 193: Cereza / Lapins / Primera ---
$ws.Range("K193").Value = "Lapins"
$ws.Range("L193").Value = "Primera"
$ws.Range("M193").Value = 630
$ws.Range("N193").Value = 4500
$ws.Range("O193").Value = 5000
$ws.Range("P193").Value = 4778
$ws.Range("Q193").Value = "$/bandeja 10 kilos"
$ws.Range("R193").Value = "Región de O'Higgins"
$ws.Range("S193").Value = 478
$ws.Range("T193").Value = 10

# --- Row 194: Cereza / Rainier / Primera ---
$ws.Range("K194").Value = "Rainier"
$ws.Range("L194").Value = "Primera"
$ws.Range("M194").Value = 480
$ws.Range("N194").Value = 10000
$ws.Range("O194").Value = 10000
$ws.Range("P194").Value = 10000
$ws.Range("Q194").Value = "$/caja 18 kilos"
$ws.Range("R194").Value = "Provincia de Curicó"
$ws.Range("S194").Value = 556
$ws.Range("T194").Value = 18

# --- Row 195: Cereza / Santina / Primera ---
$ws.Range("K195").Value = "Santina"
$ws.Range("L195").Value = "Primera"
$ws.Range("M195").Value = 450
$ws.Range("N195").Value = 5000
$ws.Range("O195").Value = 5000
$ws.Range("P195").Value = 5000
$ws.Range("Q195").Value = "$/bandeja 10 kilos"
$ws.Range("R195").Value = "Provincia de Curicó"
$ws.Range("S195").Value = 500
$ws.Range("T195").Value = 10
